# Insert a new record row at row 244 (Femacal de La Calera - Ají price list),
# shifting the existing rows 244:305 down to 245:306.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("244:244").Insert()

$ws.Range("A244").Value = 3
$ws.Range("B244").Value = "Femacal de La Calera"
$ws.Range("C244").Value = "Coquimbo"
$ws.Range("D244").Value = 44511
$ws.Range("E244").Value = 5
$ws.Range("F244").Value = 100112021
$ws.Range("G244").Value = "Ají"
$ws.Range("H244").Value = "Americana (o)"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 70
$ws.Range("K244").Value = 32000
$ws.Range("L244").Value = 33000
$ws.Range("M244").Value = 32500
$ws.Range("N244").Value = "$/caja 15 kilos"
$ws.Range("O244").Value = "Limache"
$ws.Range("P244").Value = 2167
$ws.Range("Q244").Value = 15
$ws.Range("R244").Value = "Hortaliza"
